# Applies the weekly CompStat data refresh described in the commit
# "New crime data collected": updated report title/date strings, one
# column width tweak, and refreshed crime statistics for rows 14-30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text (rich-text shared strings) ---------------------------
$ws.Range("A8").Value = "Volume 32   Number  24"
$ws.Range("C9").Value = "Report Covering the Week  6/9/2025  Through  6/15/2025"

# --- Column width tweak -------------------------------------------------
# (closest achievable value; the host's ColumnWidth setter quantizes to
# steps of 1/7 character width, so 7.433768 cannot be hit exactly)
$ws.Range("E1").ColumnWidth = 6.71

# --- Row 14 ---------------------------------------------------------
$ws.Range("N14").Value = -40

# --- Row 15 ---------------------------------------------------------
$ws.Range("F15").Value = 2
$ws.Range("I15").Value = 17
$ws.Range("K15").Value = 88.888888888888
$ws.Range("L15").Value = 183.333333333333
$ws.Range("M15").Value = 88.888888888888
$ws.Range("N15").Value = 54.545454545454

# --- Row 16 ---------------------------------------------------------
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -33.333333333333
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 61.538461538461
$ws.Range("I16").Value = 108
$ws.Range("J16").Value = 99
$ws.Range("K16").Value = 9.090909090909
$ws.Range("L16").Value = 27.058823529411
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = -46.798029556650

# --- Row 17 ---------------------------------------------------------
$ws.Range("C17").Value = 7
$ws.Range("E17").Value = -12.5
$ws.Range("F17").Value = 42
$ws.Range("G17").Value = 33
$ws.Range("H17").Value = 27.272727272727
$ws.Range("I17").Value = 174
$ws.Range("J17").Value = 157
$ws.Range("K17").Value = 10.828025477707
$ws.Range("L17").Value = 32.824427480916
$ws.Range("M17").Value = 102.325581395349
$ws.Range("N17").Value = 65.714285714285

# --- Row 18 (C18 becomes the text placeholder "0") -------------------
$ws.Range("Z1").Formula = "=CHAR(48)"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4163) | Out-Null
$ws.Range("A18").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Clear() | Out-Null

$ws.Range("D18").Value = 1
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -40
$ws.Range("J18").Value = 65
$ws.Range("K18").Value = 3.076923076923
$ws.Range("L18").Value = -4.285714285714
$ws.Range("M18").Value = -53.146853146853
$ws.Range("N18").Value = -83.037974683544

# --- Row 19 ---------------------------------------------------------
$ws.Range("C19").Value = 22
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 69.230769230769
$ws.Range("F19").Value = 64
$ws.Range("G19").Value = 60
$ws.Range("H19").Value = 6.666666666666
$ws.Range("I19").Value = 327
$ws.Range("J19").Value = 302
$ws.Range("K19").Value = 8.278145695364
$ws.Range("L19").Value = 6.514657980456
$ws.Range("M19").Value = 79.670329670329
$ws.Range("N19").Value = 73.015873015873

# --- Row 20 ---------------------------------------------------------
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 46
$ws.Range("I20").Value = 243
$ws.Range("J20").Value = 177
$ws.Range("K20").Value = 37.288135593220
$ws.Range("L20").Value = -5.813953488372
$ws.Range("M20").Value = 203.75
$ws.Range("N20").Value = -71.972318339100

# --- Row 21 ---------------------------------------------------------
$ws.Range("C21").Value = 42
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = 16.666666666666
$ws.Range("F21").Value = 181
$ws.Range("G21").Value = 143
$ws.Range("H21").Value = 26.573426573426
$ws.Range("I21").Value = 939
$ws.Range("J21").Value = 811
$ws.Range("K21").Value = 15.782983970406
$ws.Range("L21").Value = 9.440559440559
$ws.Range("M21").Value = 54.440789473684
$ws.Range("N21").Value = -47.098591549295

# --- Row 23 (D23, E23 switch from text placeholders to real numbers) --
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 250
$ws.Range("I23").Value = 38
$ws.Range("J23").Value = 27
$ws.Range("K23").Value = 40.740740740740
$ws.Range("L23").Value = 35.714285714285
$ws.Range("M23").Value = 123.529411764706

# Restore number formatting/style for D23 and E23 (they were text cells)
$ws.Range("C23").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4122) | Out-Null
$ws.Range("H23").Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4122) | Out-Null
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 100

# --- Row 24 ---------------------------------------------------------
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = 17.857142857142
$ws.Range("F24").Value = 137
$ws.Range("G24").Value = 90
$ws.Range("H24").Value = 52.222222222222
$ws.Range("I24").Value = 667
$ws.Range("J24").Value = 620
$ws.Range("K24").Value = 7.580645161290
$ws.Range("L24").Value = 5.705229793977
$ws.Range("M24").Value = 11.725293132328

# --- Row 25 ---------------------------------------------------------
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = -13.333333333333
$ws.Range("F25").Value = 64
$ws.Range("G25").Value = 42
$ws.Range("H25").Value = 52.380952380952
$ws.Range("I25").Value = 350
$ws.Range("J25").Value = 351
$ws.Range("K25").Value = -0.284900284900
$ws.Range("L25").Value = -4.891304347826

# --- Row 26 ---------------------------------------------------------
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 53
$ws.Range("G26").Value = 43
$ws.Range("H26").Value = 23.255813953488
$ws.Range("I26").Value = 241
$ws.Range("J26").Value = 247
$ws.Range("K26").Value = -2.429149797570
$ws.Range("L26").Value = 5.240174672489
$ws.Range("M26").Value = 29.569892473118

# --- Row 27 (G27, H27 switch from real numbers to text placeholders) --
$ws.Range("F27").Value = 2
$ws.Range("I27").Value = 19
$ws.Range("K27").Value = 35.714285714285
$ws.Range("L27").Value = 46.153846153846

$ws.Range("Z1").Formula = "=CHAR(48)"
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("G27").PasteSpecial(-4163) | Out-Null
$ws.Range("D27").Copy() | Out-Null
$ws.Range("G27").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Clear() | Out-Null

$ws.Range("Z1").Formula = "=""***.*"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("H27").PasteSpecial(-4163) | Out-Null
$ws.Range("E27").Copy() | Out-Null
$ws.Range("H27").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Clear() | Out-Null

# --- Row 28 ---------------------------------------------------------
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 200
$ws.Range("F28").Value = 6
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 33
$ws.Range("J28").Value = 36
$ws.Range("K28").Value = -8.333333333333
$ws.Range("L28").Value = 73.684210526315

# --- Rows 29-30 -------------------------------------------------------
$ws.Range("N29").Value = -50
$ws.Range("N30").Value = -50
